$d = $word.ActiveDocument

# 1. Insert "можно" after "глубоко" (before the comma) in the sentence.
$d.Content.Find.Execute(
    "Насколько глубоко, технически сложно",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Насколько глубоко можно, технически сложно", 2)

# 2. Remove the trailing empty paragraph at the end of the document.
$paras = $d.Paragraphs
$lastPara = $paras.Item($paras.Count)
if ($lastPara.Range.Text -eq "`r") {
    $lastPara.Range.Delete()
}
